$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.399.47'
$ws.Range('E2').Value = '  -0.54%  '
$ws.Range('D3').Value = '3.359.09'
$ws.Range('E3').Value = '  -1.82%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '571.98'
$ws.Range('E5').Value = '  -0.90%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '136.64'
$ws.Range('E6').Value = '  -1.78%  '
$ws.Range('E7').Value = '  +0.11%  '
$ws.Range('D8').Value = '3.355.62'
$ws.Range('E8').Value = '  -1.82%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.470'
$ws.Range('E9').Value = '  -1.45%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.48'
$ws.Range('E10').Value = '  -2.82%  '
$ws.Range('E11').Value = '  -2.14%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.388'
$ws.Range('E12').Value = '  -1.71%  '
$ws.Range('D13').Value = '3.934.77'
$ws.Range('E13').Value = '  -1.51%  '
$ws.Range('E14').Value = '  +1.77%  '
$ws.Range('E15').Value = '  -3.26%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '25.87'
$ws.Range('E16').Value = '  +1.11%  '
$ws.Range('D17').Value = '3.361.40'
$ws.Range('E17').Value = '  -1.24%  '
$ws.Range('D18').Value = '61.514.83'
$ws.Range('E18').Value = '  -0.46%  '
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.93'
$ws.Range('E19').Value = '  -1.18%  '
$ws.Range('B20').Value = 'Polkadot'
$ws.Range('C20').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.87'
$ws.Range('E20').Value = '  -1.24%  '
$ws.Range('E21').Value = '  -1.86%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '375.84'
$ws.Range('E22').Value = '  -3.70%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.551'
$ws.Range('E23').Value = '  -4.11%  '
$ws.Range('D24').Value = '3.508.76'
$ws.Range('E24').Value = '  -1.28%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.00'
$ws.Range('E25').Value = '  -0.05%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '71.01'
$ws.Range('E26').Value = '  -0.10%  '
$ws.Range('E27').Value = '  -2.52%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.78'
$ws.Range('E28').Value = '  +7.23%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.50'
$ws.Range('E29').Value = '  -4.16%  '
$ws.Range('E30').Value = '  +1.03%  '
$ws.Range('E31').Value = '  +3.50%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.17'
$ws.Range('E33').Value = '  -0.86%  '
$ws.Range('E34').Value = '  +0.04%  '
$ws.Range('E35').Value = '  -0.78%  '
$ws.Range('E36').Value = '  -6.18%  '
$ws.Range('E37').Value = '  -3.65%  '
$ws.Range('E38').Value = '  -2.11%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '164.80'
$ws.Range('E39').Value = '  +1.98%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0765'
$ws.Range('E40').Value = '  -4.12%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').Value = '  +0.13%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.72'
$ws.Range('E42').Value = '  -0.69%  '
$ws.Range('E43').Value = '  -0.56%  '
$ws.Range('E44').Value = '  -1.51%  '
$ws.Range('E46').Value = '  -2.59%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '23.93'
$ws.Range('E47').Value = '  +1.88%  '
$ws.Range('E48').Value = '  -2.86%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '22.84'
$ws.Range('E49').Value = '  -0.39%  '
$ws.Range('D50').Value = '2.369.24'
$ws.Range('E50').Value = '  +0.51%  '
$ws.Range('E51').Value = '  -2.19%  '
